# Updates cryptos list (price / 1h volume) table, matching the
# 'Updated cryptos list ... with GitHub Actions' commit.
#
# D/E columns hold plain text (e.g. "30.241.03", "  +6.34%  ") rather
# than numbers, so each target cell is explicitly formatted as Text
# before the new value is written -- this stops Excel's COM layer from
# auto-coercing numeric-looking strings ("1.0000", "0.9994", ...) into
# real numbers, which would silently drop meaningful trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.241.03'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +6.34%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.920.01'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.91%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.0000'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '330.30'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +5.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9994'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.78%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5227'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +3.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4084'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +4.54%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08533'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.57%  '

$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.128'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.42%  '

$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.91'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.38%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.81'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +12.58%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.450'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +4.33%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.912.03'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.55%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.415'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.62%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9993'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.78%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '95.19'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +4.56%  '

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.58%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06695'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.26%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.42'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +4.79%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.64%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.017'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.270.59'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +6.29%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.35'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.77%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.222'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.135.04'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.64%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.99'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.91%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.15'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.29%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.429'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.68%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '129.03'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.27%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +5.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1065'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.86%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.045'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +5.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.619'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.32%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02492'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.07%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06580'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.36%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2211'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.73%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.232'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +4.71%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.185'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.54%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.870'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.04%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6548'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.28%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.75'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +6.04%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.243'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.98%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6165'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.21%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.32'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.79%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.754'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.12%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.090'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +4.91%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.248'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.74%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.27'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.92%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.164'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.10%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.69'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +5.07%  '
